$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row data (rows 42-45)
$rows = @(
    @{ Row=42; A=45784.4906693287;  B="4G_PTO015M_HNI"; C="POWER_AC_EAS"; D="Thành công"; E="Phúc Thọ";   F="GS vận hành MFĐ PTO-KTĐH" },
    @{ Row=43; A=45784.49078851852; B="3G_BVI110M_HNI"; C="SITE_OOS";     D="Thành công"; E="Ba Vì";      F="GS vận hành MFĐ BVI-KTĐH" },
    @{ Row=44; A=45784.49083607639; B="4G-BVI110M-HNI"; C="SITE_OOS";     D="Thành công"; E="Ba Vì";      F="GS vận hành MFĐ BVI-KTĐH" },
    @{ Row=45; A=45784.54312510684; B="SR_TTT032M_HNI"; C="POWER_AC_EAS"; D="Thành công"; E="Thạch Thất"; F="GS vận hành MFĐ-BTS -TTT-KTĐH" }
)

# Write the brand-new unique string values first, in the exact order they
# must be appended to the shared string table: the three new "Ten NE"
# values, then the new district value, then the new Zalo-group value.
$ws.Cells.Item(42, 2).Value = "4G_PTO015M_HNI"
$ws.Cells.Item(43, 2).Value = "3G_BVI110M_HNI"
$ws.Cells.Item(44, 2).Value = "4G-BVI110M-HNI"
$ws.Cells.Item(42, 5).Value = "Phúc Thọ"
$ws.Cells.Item(42, 6).Value = "GS vận hành MFĐ PTO-KTĐH"

# Now fill in the remaining cells (these all reuse already-existing shared
# strings, so the order here does not affect the shared string table).
foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    $ws.Cells.Item($rowIndex, 5).Value = $r.E
    $ws.Cells.Item($rowIndex, 6).Value = $r.F
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
}
